# This script re-applies a weekly refresh of the "Chirimoya" price sheet.
# The rows' Fecha/Calidad/Volumen/Precio.../Unidad/Origen data are shuffled
# amongst themselves (each row ends up showing the data that another row in
# the original sheet had), while all the other columns (A,B,C,E,F,G,H,I,J,K)
# stay untouched.
#
# Mapping: new row r gets the D/L/M/N/O/P/Q/R/S/T values that used to live
# on row perm[r] in the original sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$perm = @{
    2  = 6
    3  = 18
    4  = 17
    5  = 9
    6  = 7
    7  = 16
    8  = 10
    9  = 3
    10 = 4
    11 = 13
    12 = 8
    13 = 19
    14 = 2
    15 = 5
    16 = 11
    17 = 15
    18 = 12
    19 = 20
    20 = 14
}

$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# First, snapshot the "before" values for every touched cell, because the
# permutation reads from rows that will themselves be overwritten later in
# the loop.
$snapshot = @{}
foreach ($r in $perm.Keys) {
    foreach ($col in $cols) {
        $addr = "$col$r"
        $snapshot[$addr] = $ws.Range($addr).Value2
    }
}

# Now apply: for each destination row, write the snapshot values captured
# from its source row.
foreach ($r in $perm.Keys) {
    $srcRow = $perm[$r]
    foreach ($col in $cols) {
        $srcAddr = "$col$srcRow"
        $dstAddr = "$col$r"
        $ws.Range($dstAddr).Value2 = $snapshot[$srcAddr]
    }
}
